# v1.2 - An error message will be displayed if the article body exceeds 1000 words
# Adds a new assumption row (LH-SRS-ASSUMP-002) to LH_SRS_ASSUMPTIONS and a matching
# version-history entry (V1.2) to "Version Histroy".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH_SRS_ASSUMPTIONS")
$ws2 = $wb.Worksheets.Item("Version Histroy")

# ---------------------------------------------------------------------------
# LH_SRS_ASSUMPTIONS sheet
# ---------------------------------------------------------------------------

# Row 2 picked up a horizontal "align left" on top of its existing
# vertical-center / wrap formatting (date cell keeps its own numeric format).
$ws1.Range("A2").HorizontalAlignment = -4131
$ws1.Range("A2").NumberFormat = "mm-dd-yy"
$ws1.Range("B2:J2").HorizontalAlignment = -4131

# Shrink row 2 now that the content fits in less height.
$ws1.Rows.Item(2).RowHeight = 90

# New assumption row.
$ws1.Range("A3").Value = 45773
$ws1.Range("A3").HorizontalAlignment = -4131
$ws1.Range("A3").NumberFormat = "mm-dd-yy"

$ws1.Range("B3").Value = "LH-SRS-ASSUMP-002"
$ws1.Range("C3").Value = "SRS-PUB-004"
$ws1.Range("D3").Value = "An error message will be displayed if the article body exceeds 1000 words"
$ws1.Range("E3").Value = "The system will not show an error message if the article exceeds 1000 words, and the user will not be informed that the article cannot be submitted due to the word limit being exceeded"
$ws1.Range("F3").Value = "Hala Eldaly"
$ws1.Range("G3").Value = "High"
$ws1.Range("H3").Value = "Pending"
$ws1.Range("J3").Value = "I made an assumption on the error message that will be shown if the user tries to submit an article that exceeds 1000 words: [`"Exceeds 1000 words`"].`n"

$ws1.Range("B3:J3").HorizontalAlignment = -4131

$ws1.Rows.Item(3).RowHeight = 117.75

# Column J needed to widen to fit the new note text.
$ws1.Columns.Item(10).ColumnWidth = 29.3

# ---------------------------------------------------------------------------
# Version Histroy sheet
# ---------------------------------------------------------------------------

$ws2.Range("A4").Value = "V1.2"
$ws2.Range("B4").Value = "Hala  Eldaly"
$ws2.Range("C4").Value = "LH-SRS-ASSUMP-002"
$ws2.Range("D4").Value = 45773

# ---------------------------------------------------------------------------
# Selections (restore LH_SRS_ASSUMPTIONS as the active/visible sheet).
# ---------------------------------------------------------------------------

$ws2.Range("D13").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B8").Select() | Out-Null
